# Update the HITS@k tables so the "supporting entity" visualization
# differentiates numerical vs. non-numerical answers.
#
# For every changed row we:
#   1. Write the new raw hit count into column B.
#   2. Force column D to remain a plain text percentage string (rather
#      than letting Excel auto-convert "NN.NN%" into a numeric percent
#      cell) by flipping the cell to a text NumberFormat before the
#      write, then resetting the cell style back to "Normal" so no
#      stray style index is left attached to the cell.
$wb = $excel.ActiveWorkbook

function Set-HitRow {
    param(
        [object]$ws,
        [int]$row,
        [int]$hits,
        [string]$pct
    )
    $ws.Cells.Item($row, 2).Value = $hits
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $pct
    $dCell.Style = "Normal"
}

# --- Total Hits ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Total Hits")
Set-HitRow $ws 2 1098 "37.71%"
Set-HitRow $ws 3 2171 "37.28%"
Set-HitRow $ws 4 3262 "37.34%"
Set-HitRow $ws 5 4371 "37.53%"
Set-HitRow $ws 6 5458 "37.49%"

# --- Hits_entity ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Hits_entity")
Set-HitRow $ws 3 917  "32.49%"
Set-HitRow $ws 4 1365 "32.25%"
Set-HitRow $ws 5 1834 "32.49%"
Set-HitRow $ws 6 2294 "32.52%"

# --- Hits_numerical --------------------------------------------------------
$ws = $wb.Worksheets.Item("Hits_numerical")
Set-HitRow $ws 2 156 "23.82%"
Set-HitRow $ws 3 322 "24.58%"
Set-HitRow $ws 4 492 "25.04%"
Set-HitRow $ws 5 658 "25.11%"
Set-HitRow $ws 6 809 "24.70%"

# --- Hits_date -------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hits_date")
Set-HitRow $ws 3 220 "41.51%"
Set-HitRow $ws 4 328 "41.26%"
Set-HitRow $ws 5 439 "41.42%"
Set-HitRow $ws 6 562 "42.42%"
